$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address, new text value.
$updates = @(
    @('D2', '29.045.24'),
    @('E2', '  -0.01%  '),
    @('D3', '1.834.69'),
    @('E3', '  +0.32%  '),
    @('D4', '0.9991'),
    @('E4', '  +0.04%  '),
    @('D5', '244.25'),
    @('E5', '  +1.45%  '),
    @('D6', '0.6335'),
    @('E6', '  +1.97%  '),
    @('E7', '  +0.10%  '),
    @('D8', '0.07586'),
    @('E9', '  +0.92%  '),
    @('D10', '22.78'),
    @('E10', '  +0.32%  '),
    @('D11', '0.07743'),
    @('E11', '  +0.92%  '),
    @('D12', '1.835.89'),
    @('E12', '  +0.41%  '),
    @('D13', '4.992'),
    @('E13', '  +0.58%  '),
    @('D14', '0.6710'),
    @('E14', '  +1.23%  '),
    @('D15', '83.21'),
    @('E15', '  +1.37%  '),
    @('D16', '0.000009857'),
    @('E16', '  +8.72%  '),
    @('D17', '6.116'),
    @('E17', '  +1.45%  '),
    @('D18', '29.079.59'),
    @('D19', '12.56'),
    @('E19', '  +1.54%  '),
    @('D20', '226.93'),
    @('E20', '  +0.62%  '),
    @('D21', '0.9996'),
    @('E21', '  -0.05%  '),
    @('D22', '7.220'),
    @('E22', '  +1.00%  '),
    @('E23', '  +0.05%  '),
    @('D24', '160.50'),
    @('E24', '  +0.57%  '),
    @('D25', '0.1405'),
    @('E25', '  +3.60%  '),
    @('D26', '8.541'),
    @('E26', '  +1.41%  '),
    @('D27', '17.96'),
    @('E27', '  +0.87%  '),
    @('E28', '  +0.50%  '),
    @('D29', '4.119'),
    @('E29', '  +1.59%  '),
    @('D30', '4.055'),
    @('E30', '  +0.41%  '),
    @('E31', '  +0.32%  '),
    @('D32', '0.05403'),
    @('E32', '  +3.05%  '),
    @('D33', '1.860'),
    @('E33', '  +0.78%  '),
    @('D34', '0.7469'),
    @('E34', '  +2.02%  '),
    @('D35', '1.142'),
    @('E35', '  -0.94%  '),
    @('D36', '2.672'),
    @('E36', '  +0.95%  '),
    @('D37', '1.243.96'),
    @('E37', '  -3.55%  '),
    @('E38', '  +0.78%  '),
    @('E39', '  +0.39%  '),
    @('D40', '6.619'),
    @('E40', '  +5.03%  '),
    @('D41', '0.9049'),
    @('E41', '  +0.39%  '),
    @('E42', '  +0.21%  '),
    @('D43', '102.78'),
    @('E43', '  +0.88%  '),
    @('D44', '1.986.66'),
    @('E44', '  +0.56%  '),
    @('D45', '0.00000000123'),
    @('E45', '  +2.88%  '),
    @('D46', '64.77'),
    @('E46', '  +1.32%  '),
    @('E47', '  -0.02%  '),
    @('D48', '0.4106'),
    @('E48', '  +3.53%  '),
    @('D49', '9.057'),
    @('E49', '  +2.73%  '),
    @('D50', '0.05785'),
    @('E50', '  +0.06%  '),
    @('D51', '6.773'),
    @('E51', '  +1.57%  ')
)

foreach ($u in $updates) {
    $addr = $u[0]
    $text = $u[1]
    $cell = $ws.Range($addr)
    # Numbers-like strings (e.g. "0.9991") get auto-converted to the
    # numeric type by Excel on plain assignment; force text storage
    # via a temporary Text number format, then restore the default
    # (General/Normal) style so no visible formatting change remains.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}
